$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (the "PV 002" shared string is renamed to "PV 001" via B2,
# and row quantities are updated)
$ws.Range("B2").Value = "PV 001"
$ws.Range("C2").Value = 10
$ws.Range("C3").Value = 30
$ws.Range("C4").Value = 50

# Update the selected cell in the sheet view
$ws.Range("B9").Select()
